# Atualiza notas dos alunos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Filipe Palma Abreu" was dropped from the roster: delete his row (row 3),
# shifting "Gabriel Andrade Vieira" and "Lucas Borges Jagersbacher" up.
$ws.Range("A3").EntireRow.Delete()

# New "Total" / "Conceito" columns, styled like the other header cells (bold).
$ws.Range("G1").Value = "Total"
$ws.Range("H1").Value = "Conceito"
$ws.Range("G1:H1").Font.Bold = $true
$ws.Range("G1").ColumnWidth = 7.45

# Grades entered for each student (C1/C2 columns).
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 1.5
$ws.Range("C3").Value = 1.5

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2

# Totals per student.
$ws.Range("G2").Formula = "=SUM(B2:F2)"
$ws.Range("G3:G4").Formula = "=SUM(B3:F3)"

[void]$ws.Range("A7").Select()
